# Updated cryptos list on Sun Jan 28 03:06:15 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '42.111.63'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.60%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.280.11'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +0.19%  '

$ws.Range("E4").Value = '  -0.03%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '155.32'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +15,408.24%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '305.21'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.29%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '94.32'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +1.26%  '

$ws.Range("E8").Value = '  -0.36%  '

$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("E10").Value = '  +0.80%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '34.17'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +4.71%  '

$ws.Range("E12").Value = '  -0.05%  '

$ws.Range("E13").Value = '  -2.24%  '

$ws.Range("E14").Value = '  -0.50%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '2.632.57'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.16%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '14.41'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +0.87%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.279.40'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -1.10%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.794'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +3.89%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '42.032.73'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.56%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '12.74'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +3.60%  '

$ws.Range("E21").Value = '  +0.72%  '

$ws.Range("E22").Value = '  +0.69%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '68.01'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +0.90%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '243.78'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.01%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.59'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.44%  '

$ws.Range("E26").Value = '  +0.25%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.13%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '24.08'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.81%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '35.95'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +5.25%  '

$ws.Range("E31").Value = '  +1.33%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '160.82'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +1.15%  '

$ws.Range("E33").Value = '  +2.82%  '

$ws.Range("E34").Value = '  +0.02%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.0753'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +0.13%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '3.09'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +1.09%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '17.03'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +1.87%  '

$ws.Range("E39").Value = '  -0.24%  '

$ws.Range("E40").Value = '  -0.25%  '

$ws.Range("E41").Value = '  -0.97%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '4.23'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +7.23%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '2.018.97'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.55%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '19.73'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.30%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.27'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +11.27%  '

$ws.Range("E46").Value = '  +1.51%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '10.23'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -2.22%  '

$ws.Range("E48").Value = '  -0.20%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '53.54'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +3.35%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.51'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.95%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '72.21'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -1.16%  '
